# Updated cryptos list values (Price / Volume(1h)) per target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay text even when the string looks numeric,
    # matching the existing inline-string cells in this sheet, then
    # drop back to the default (unstyled) cell so no style is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "62.821.35"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "3.471.04"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue $ws.Range("D5") "413.91"
$ws.Range("E5").Value = "  +1.13%  "
Set-TextValue $ws.Range("D6") "130.71"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.22%  "
Set-TextValue $ws.Range("D10") "0.147"
$ws.Range("E10").Value = "  +3.76%  "
Set-TextValue $ws.Range("D11") "42.68"
$ws.Range("E11").Value = "  -0.40%  "
Set-TextValue $ws.Range("D12") "9.58"
$ws.Range("E12").Value = "  +3.35%  "
Set-TextValue $ws.Range("D13") "0.0000220"
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("D14").Value = "4.021.45"
$ws.Range("E14").Value = "  +1.19%  "
Set-TextValue $ws.Range("D15") "0.141"
$ws.Range("E15").Value = "  -0.21%  "
Set-TextValue $ws.Range("D16") "20.53"
$ws.Range("E16").Value = "  -4.12%  "
$ws.Range("D17").Value = "3.463.91"
$ws.Range("E17").Value = "  +1.77%  "
Set-TextValue $ws.Range("D18") "12.65"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "62.732.39"
$ws.Range("E20").Value = "  +1.06%  "
Set-TextValue $ws.Range("D21") "466.09"
$ws.Range("E21").Value = "  +2.73%  "
Set-TextValue $ws.Range("D22") "90.75"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +0.98%  "
Set-TextValue $ws.Range("D25") "10.61"
$ws.Range("E25").Value = "  +16.46%  "
Set-TextValue $ws.Range("D26") "3.30"
$ws.Range("E26").Value = "  +1.42%  "
Set-TextValue $ws.Range("D27") "33.42"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  -1.63%  "
Set-TextValue $ws.Range("D34") "40.81"
$ws.Range("E34").Value = "  -4.94%  "
$ws.Range("E35").Value = "  +0.03%  "
Set-TextValue $ws.Range("D36") "58.75"
$ws.Range("E36").Value = "  +7.99%  "
Set-TextValue $ws.Range("D37") "0.0491"
$ws.Range("E37").Value = "  -2.50%  "
Set-TextValue $ws.Range("D38") "0.999"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +4.06%  "
Set-TextValue $ws.Range("D40") "147.29"
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("E44").Value = "  +7.24%  "
$ws.Range("E45").Value = "  +3.89%  "
Set-TextValue $ws.Range("D46") "4.34"
$ws.Range("E46").Value = "  +1.79%  "
Set-TextValue $ws.Range("D47") "2.41"
$ws.Range("E47").Value = "  +12.94%  "
$ws.Range("D48").Value = "0.0₃0556"
$ws.Range("E48").Value = "  +26.83%  "
$ws.Range("E49").Value = "  -1.82%  "
Set-TextValue $ws.Range("D50") "22.33"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  +0.90%  "
